# Generate Report for Handback
# Updates the handback-status workbook with the freshly generated report
# values: new source-file UUIDs/paths and new handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

# ---- new identifiers -------------------------------------------------
$uuid1New = "a9ca9bc8-6ed6-41bd-bc60-9c41a806a04e"
$uuid2New = "ffff2e8876e5-c83c-47f2-97b6-5dd736ee3d18"
$hashNew  = "484e9777cdc6ef738bd963f8dd4327dba3002178"

$xlfZhCn = "$uuid1New.$hashNew.zh-cn.xlf"
$xlfDeDe = "$uuid1New.$hashNew.de-de.xlf"

$dateGenerate = "2016-08-28 15:02:33"
$dateZhCnHandoff  = "2016-08-28 15:02:28"
$dateZhCnHandback = "2016-08-28 15:02:44"
$dateDeDeHandoff  = "2016-08-28 15:02:33"
$dateDeDeHandback = "2016-08-28 15:02:51"

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$uuid1New.md"
$wsOverview.Range("B2").Value = "e2e\$uuid1New.md"
$wsOverview.Range("G2").Value = $dateGenerate

$wsOverview.Range("A3").Value = "$uuid2New.md"
$wsOverview.Range("B3").Value = "e2e\$uuid2New.md"
$wsOverview.Range("G3").Value = $dateGenerate

$i = 0
foreach ($hl in $wsOverview.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $hl.TextToDisplay = "e2e\$uuid1New.md"
    } else {
        $hl.TextToDisplay = "e2e\$uuid2New.md"
    }
}

# ---- zh-cn sheet ---------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$uuid1New.md"
$wsZhCn.Range("I2").Value = "$uuid1New.md"
$wsZhCn.Range("G2").Value = $xlfZhCn
$wsZhCn.Range("J2").Value = $xlfZhCn
$wsZhCn.Range("H2").Value = $dateZhCnHandoff
$wsZhCn.Range("K2").Value = $dateZhCnHandback

$wsZhCn.Range("A3").Value = "$uuid2New.md"
$wsZhCn.Range("I3").Value = "$uuid2New.md"
$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("J3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = $dateZhCnHandoff
$wsZhCn.Range("K3").Value = $dateZhCnHandback

$i = 0
foreach ($hl in $wsZhCn.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $hl.TextToDisplay = "$uuid1New.md"
    } elseif ($i -eq 2) {
        $hl.TextToDisplay = "$uuid1New.md"
    } elseif ($i -eq 3) {
        $hl.TextToDisplay = "$uuid2New.md"
    } else {
        $hl.TextToDisplay = "$uuid2New.md"
    }
}

# ---- de-de sheet ---------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$uuid1New.md"
$wsDeDe.Range("I2").Value = "$uuid1New.md"
$wsDeDe.Range("G2").Value = $xlfDeDe
$wsDeDe.Range("J2").Value = $xlfDeDe
$wsDeDe.Range("H2").Value = $dateDeDeHandoff
$wsDeDe.Range("K2").Value = $dateDeDeHandback

$wsDeDe.Range("A3").Value = "$uuid2New.md"
$wsDeDe.Range("I3").Value = "$uuid2New.md"
$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("J3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = $dateDeDeHandoff
$wsDeDe.Range("K3").Value = $dateDeDeHandback

$i = 0
foreach ($hl in $wsDeDe.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $hl.TextToDisplay = "$uuid1New.md"
    } elseif ($i -eq 2) {
        $hl.TextToDisplay = "$uuid1New.md"
    } elseif ($i -eq 3) {
        $hl.TextToDisplay = "$uuid2New.md"
    } else {
        $hl.TextToDisplay = "$uuid2New.md"
    }
}
